$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.628.86'
$ws.Range('D3').Value = '2.332.14'
$ws.Range('E3').Value = '  +2.02%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Formula = '="546.61"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +2.10%  '
$ws.Range('D6').Formula = '="131.29"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Formula = '="0.580"'
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  -0.92%  '
$ws.Range('D9').Value = '2.331.10'
$ws.Range('E9').Value = '  +2.21%  '
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('D11').Formula = '="5.50"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  +0.60%  '
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').Formula = '="0.336"'
$ws.Range('D13').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').Formula = '="23.81"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '60.580.39'
$ws.Range('E15').Value = '  +4.16%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '2.748.02'
$ws.Range('E16').Value = '  +2.03%  '
$ws.Range('D17').Formula = '="0.0000133"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('D18').Value = '2.342.03'
$ws.Range('E18').Value = '  +1.85%  '
$ws.Range('D19').Formula = '="10.63"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  +0.90%  '
$ws.Range('D20').Formula = '="4.16"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  -0.56%  '
$ws.Range('D21').Formula = '="316.02"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('E22').Value = '  +1.78%  '
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('D24').Formula = '="64.22"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +1.61%  '
$ws.Range('D25').Formula = '="0.171"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +2.09%  '
$ws.Range('D26').Formula = '="0.999"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('D27').Formula = '="7.86"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  -1.45%  '
$ws.Range('D28').Formula = '="1.34"'
$ws.Range('D28').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  +4.26%  '
$ws.Range('D29').Formula = '="1.20"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  +11.09%  '
$ws.Range('D30').Formula = '="173.12"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  +1.51%  '
$ws.Range('D31').Formula = '="1.73"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  +2.06%  '
$ws.Range('D32').Value = '0.0₃0734'
$ws.Range('E32').Value = '  +1.49%  '
$ws.Range('D33').Formula = '="5.95"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  +2.57%  '
$ws.Range('D34').Formula = '="1.39"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +12.29%  '
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('B36').Value = 'USDe'
$ws.Range('C36').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D36').Formula = '="0.999"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Formula = '="17.90"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').Formula = '="4.08"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  +4.36%  '
$ws.Range('D40').Formula = '="328.73"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +13.75%  '
$ws.Range('E41').Value = '  +2.80%  '
$ws.Range('D42').Formula = '="37.99"'
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('D43').Formula = '="138.49"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('D44').Formula = '="3.49"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +1.65%  '
$ws.Range('D45').Formula = '="0.0944"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('D46').Formula = '="19.39"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +7.08%  '
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('E49').Value = '  +1.70%  '
$ws.Range('D50').Value = '0.0₆0217'
$ws.Range('E50').Value = '  +18.47%  '
$ws.Range('D51').Formula = '="11.03"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +0.71%  '
$excel.CutCopyMode = 0

